# dados_ibge_cidade.xlsx — replace the placeholder "Pequena/Media/Grande" +
# "Municipio" template with real IBGE data for five Rondonia municipalities,
# and drop the now-unused leading classification column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Header row (becomes columns A:J once the old "Pequena, Media ou
#    Grande ?" column is gone and everything else shifts one to the left).
# ---------------------------------------------------------------------
$headers = @(
    "Municípios",
    "UF",
    "Cod. Municipio",
    "População no último censo",
    "Salário médio mensal dos trabalhadores formais",
    "Matrículas no ensino fundamental",
    "PIB per capita",
    "Mortalidade Infantil",
    "Área urbanizada",
    "Área da unidade territorial"
)
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}

# ---------------------------------------------------------------------
# 2. Data rows 2..6 — five municipalities from Rondônia (RO).
#    Column D ("População no último censo") is left blank in every row,
#    same as in the source sheet.
# ---------------------------------------------------------------------
$data = @(
    @("Alta Floresta D'Oeste", "RO", "1100015", "", "7.067,127", "1,8",    "3.051",  "32.619,88", "5,43",  "6,46"),
    @("Ariquemes",             "RO", "1100023", "", "4.426,571", "2,1",    "13.078", "28.878,27", "21,26", "33,26"),
    @("Cabixi",                "RO", "1100031", "", "1.314,352", "2",      "729",    "47.051,83", "-",     "2,24"),
    @("Cacoal",                "RO", "1100049", "", "3.793,000", "1,9",    "10.969", "32.313,26", "9,61",  "28,04"),
    @("Cerejeiras",            "RO", "1100056", "", "2.783,300", "2,1",    "2.171",  "46.185,81", "3,68",  "6,31")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        if ($c -ne 3) {
            $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
        }
    }
}

# Column D stays empty but keeps the "0.00" number format used by the
# original template for that column.
$ws.Range("D2:D8").NumberFormat = "0.00"

# ---------------------------------------------------------------------
# 3. Header formatting: columns G,H,I ("PIB per capita",
#    "Mortalidade Infantil", "Área urbanizada") use the highlighted
#    header style; everything else (incl. the new last column J) uses
#    the plain bordered header style.
# ---------------------------------------------------------------------
$ws.Range("H1").Copy()
$ws.Range("G1:I1").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 4. Drop the old trailing columns (K was "Área da unidade territorial"
#    which is now column J; L/M were always unused) and the two now
#    unused bottom rows, so the sheet's used range shrinks to A1:J8.
# ---------------------------------------------------------------------
$ws.Range("K1:M10").Clear()
$ws.Columns.Item(11).ColumnWidth = 8.43

$ws.Rows.Item(9).Delete()
$ws.Rows.Item(9).Delete()

# ---------------------------------------------------------------------
# 5. Resize the table to the new 10-column, 6-row extent (header +
#    5 data rows). Column names are re-read from row 1 automatically.
# ---------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:J6"))

# ---------------------------------------------------------------------
# 6. Column widths, matching the final layout (values are the XML
#    "width" units minus the fixed 0.8333... padding this host adds).
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth  = 27.830729166666668
$ws.Columns.Item(2).ColumnWidth  = 27.830729166666668
$ws.Columns.Item(3).ColumnWidth  = 24.166666666666668
$ws.Columns.Item(4).ColumnWidth  = 47.498697916666664
$ws.Columns.Item(5).ColumnWidth  = 56.608072916666664
$ws.Columns.Item(6).ColumnWidth  = 41.721354166666664
$ws.Columns.Item(7).ColumnWidth  = 19.830729166666668
$ws.Columns.Item(8).ColumnWidth  = 22.608072916666668
$ws.Columns.Item(9).ColumnWidth  = 22.608072916666668
$ws.Columns.Item(10).ColumnWidth = 33.830729166666664

# ---------------------------------------------------------------------
# 7. Selection, matching where the author last clicked.
# ---------------------------------------------------------------------
$ws.Range("L8").Select()
